$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '307.00'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '-4.52%'

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '39.16'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '-8.09%'

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '5.099'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '-2.32%'

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.07684'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '-6.04%'

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '4.247'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '-1.67%'

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.641'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '-7.87%'

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.9150'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '-3.46%'

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.1017'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '-9.25%'

$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '-7.67%'

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.09066'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '-3.20%'

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.04423'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '-4.38%'

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.1056'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '-0.16%'

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.001260'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '-3.02%'

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.005808'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '0.72%'

$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '2,416.75%'

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '3.361'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '-0.03%'

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '2.419'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '-4.82%'

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.3308'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '-1.74%'

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '7.046'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '-5.44%'

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.1347'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '-3.00%'

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.2813'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '10.46%'

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.04136'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '-0.34%'

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.001203'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '-3.80%'

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.004100'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '-4.40%'

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.0001299'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '6.57%'

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.02437'

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.05196'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '-7.27%'

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.007912'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '-2.46%'

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.1317'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '-5.99%'

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.007165'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '9.48%'

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.001948'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '-3.96%'

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.008377'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '9.67%'

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.3336'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '-4.19%'

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.00006429'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '-5.00%'

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.00000000750'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '0.02%'

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.003001'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '-26.76%'

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.004362'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '41.98%'

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.00002099'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '0.02%'

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0001999'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '0.02%'
